# Rerun all experiment : convert linear regression to logistic regression
#
# Updates the "R-Proposed" ranking column (D) on the three per-metric sheets
# with the re-run (logistic regression) results, swaps one player in the
# proposed Top-10 ranking (Quincy Promes / Matthias Ginter) together with his
# updated market value, and refreshes the selected cell on several sheets to
# match the author's last cursor position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Top 10 players goal 90" - R-Proposed (D) column rerun values
# ---------------------------------------------------------------------
$wsGoal = $wb.Worksheets.Item("Top 10 players goal 90")
$wsGoal.Range("D2").Value = 60
$wsGoal.Range("D3").Value = 121
$wsGoal.Range("D5").Value = 106
$wsGoal.Range("D7").Value = 85
$wsGoal.Range("D8").Value = 113
$wsGoal.Range("D9").Value = 54
$wsGoal.Range("D10").Value = 79
$wsGoal.Range("D11").Value = 23

# ---------------------------------------------------------------------
# Sheet "Top 10 players assist 90" - R-Proposed (D) column rerun values
# plus a Result (E) flip on row 11, and the refreshed cursor cell.
# ---------------------------------------------------------------------
$wsAssist = $wb.Worksheets.Item("Top 10 players assist 90")
$wsAssist.Range("D2").Value = 6
$wsAssist.Range("D3").Value = 77
$wsAssist.Range("D4").Value = 106
$wsAssist.Range("D5").Value = 5
$wsAssist.Range("D6").Value = 101
$wsAssist.Range("D7").Value = 82
$wsAssist.Range("D8").Value = 58
$wsAssist.Range("D9").Value = 85
$wsAssist.Range("D10").Value = 23
$wsAssist.Range("D11").Value = 30
$wsAssist.Range("E11").Value = "Lose"
$wsAssist.Activate()
$wsAssist.Range("F15").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Top 10 players goal assist 90" - R-Proposed (D) column rerun
# values.
# ---------------------------------------------------------------------
$wsGoalAssist = $wb.Worksheets.Item("Top 10 players goal assist 90")
$wsGoalAssist.Range("D2").Value = 6
$wsGoalAssist.Range("D3").Value = 106
$wsGoalAssist.Range("D4").Value = 60
$wsGoalAssist.Range("D5").Value = 85
$wsGoalAssist.Range("D6").Value = 121
$wsGoalAssist.Range("D8").Value = 77
$wsGoalAssist.Range("D10").Value = 113
$wsGoalAssist.Range("D11").Value = 5

# ---------------------------------------------------------------------
# Sheet "Top 10 VAEP ranking" - no data changed, only the saved cursor
# position moved.
# ---------------------------------------------------------------------
$wsVaep = $wb.Worksheets.Item("Top 10 VAEP ranking")
$wsVaep.Activate()
$wsVaep.Range("G13").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Top 10 proposed ranking" - the rerun reshuffles three rows
# (Steven Berghuis / Pau Torres / Dejan Kulusevski rotate down one slot)
# and swaps the 9th/10th entries, with Quincy Promes being replaced by
# Matthias Ginter (and his market value) in the new ranking, which also
# updates the total market value footer.
# ---------------------------------------------------------------------
$wsProposed = $wb.Worksheets.Item("Top 10 proposed ranking")
$wsProposed.Range("B6").Value = "Pau Torres"
$wsProposed.Range("C6").Value = "20 million euro"
$wsProposed.Range("B7").Value = "Dejan Kulusevski"
$wsProposed.Range("C7").Value = "35 million euro"
$wsProposed.Range("B8").Value = "Steven Berghuis"
$wsProposed.Range("C8").Value = "14 million euro"
$wsProposed.Range("B10").Value = "William Carvalho"
$wsProposed.Range("C10").Value = "16 million euro"
$wsProposed.Range("B11").Value = "Matthias Ginter"
$wsProposed.Range("C11").Value = "28 million euro"
$wsProposed.Range("C13").Value = "242,3 million euro"
$wsProposed.Activate()
$wsProposed.Range("F16").Select() | Out-Null
